# Auto-generated edit script applying numeric corrections to the
# Halicarnassus_Profits leve-crafting profit tables across all job sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 670.2857
$ws.Range("I6").Value = 36.555557
$ws.Range("J6").Value = 1811
$ws.Range("K6").Value = 109.666671
$ws.Range("L6").Value = 5433
$ws.Range("M6").Value = 2.333328999999992
$ws.Range("N6").Value = -5657
$ws.Range("H15").Value = 1481.3654
$ws.Range("I15").Value = 1481.3654
$ws.Range("K15").Value = 4444.0962
$ws.Range("M15").Value = -4275.0962
$ws.Range("H29").Value = 7498.75
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 7498.75
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 22496.25
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -23058.25
$ws.Range("H58").Value = 1577.1818
$ws.Range("J58").Value = 2827
$ws.Range("L58").Value = 8481
$ws.Range("N58").Value = -8781
$ws.Range("H88").Value = 1257.7142
$ws.Range("J88").Value = 1574.75
$ws.Range("L88").Value = 1574.75
$ws.Range("N88").Value = -2386.75
$ws.Range("H91").Value = 1257.7142
$ws.Range("J91").Value = 1574.75
$ws.Range("L91").Value = 1574.75
$ws.Range("N91").Value = -4382.75
$ws.Range("H100").Value = 4860.636
$ws.Range("I100").Value = 2122
$ws.Range("J100").Value = 6425.5713
$ws.Range("K100").Value = 2122
$ws.Range("L100").Value = 6425.5713
$ws.Range("M100").Value = -1581
$ws.Range("N100").Value = -7507.5713

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 120
$ws.Range("I5").Value = 120
$ws.Range("K5").Value = 120
$ws.Range("M5").Value = -8
$ws.Range("H45").Value = 3947.4546
$ws.Range("J45").Value = 5399.4
$ws.Range("L45").Value = 5399.4
$ws.Range("N45").Value = -6153.4
$ws.Range("H88").Value = 1133.9333
$ws.Range("I88").Value = 743.4286
$ws.Range("J88").Value = 1475.625
$ws.Range("K88").Value = 743.4286
$ws.Range("L88").Value = 1475.625
$ws.Range("M88").Value = -337.4286
$ws.Range("N88").Value = -2287.625
$ws.Range("H91").Value = 1133.9333
$ws.Range("I91").Value = 743.4286
$ws.Range("J91").Value = 1475.625
$ws.Range("K91").Value = 743.4286
$ws.Range("L91").Value = 1475.625
$ws.Range("M91").Value = 660.5714
$ws.Range("N91").Value = -4283.625
$ws.Range("H101").Value = 55554.5
$ws.Range("J101").Value = 55554.5
$ws.Range("L101").Value = 55554.5
$ws.Range("N101").Value = -62044.5
$ws.Range("H122").Value = 1705.375
$ws.Range("I122").Value = 1482.1666
$ws.Range("K122").Value = 4446.4998
$ws.Range("M122").Value = -1996.4998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 120
$ws.Range("I4").Value = 120
$ws.Range("K4").Value = 120
$ws.Range("M4").Value = -5
$ws.Range("H99").Value = 8636
$ws.Range("I99").Value = 7954.5
$ws.Range("J99").Value = 9999
$ws.Range("K99").Value = 7954.5
$ws.Range("L99").Value = 9999
$ws.Range("M99").Value = -6456.5
$ws.Range("N99").Value = -12995

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 974.75
$ws.Range("I16").Value = 934
$ws.Range("J16").Value = 1749
$ws.Range("K16").Value = 934
$ws.Range("L16").Value = 1749
$ws.Range("M16").Value = -647
$ws.Range("N16").Value = -2323
$ws.Range("H35").Value = 183.33333
$ws.Range("I35").Value = 156.875
$ws.Range("J35").Value = 395
$ws.Range("K35").Value = 156.875
$ws.Range("L35").Value = 395
$ws.Range("M35").Value = 137.125
$ws.Range("N35").Value = -983
$ws.Range("H50").Value = 27600
$ws.Range("I50").Value = 19333.334
$ws.Range("K50").Value = 19333.334
$ws.Range("M50").Value = -18708.334
$ws.Range("H51").Value = 32781.25
$ws.Range("J51").Value = 39708.332
$ws.Range("L51").Value = 39708.332
$ws.Range("N51").Value = -41180.332
$ws.Range("H60").Value = 1048.25
$ws.Range("I60").Value = 1048.25
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 1048.25
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -537.25
$ws.Range("N60").ClearContents()
$ws.Range("H61").Value = 32781.25
$ws.Range("J61").Value = 39708.332
$ws.Range("L61").Value = 39708.332
$ws.Range("N61").Value = -40404.332
$ws.Range("H92").Value = 39999
$ws.Range("J92").Value = 39999
$ws.Range("L92").Value = 39999
$ws.Range("N92").Value = -44991
$ws.Range("H113").Value = 974.75
$ws.Range("I113").Value = 934
$ws.Range("J113").Value = 1749
$ws.Range("K113").Value = 934
$ws.Range("L113").Value = 1749
$ws.Range("M113").Value = 1236
$ws.Range("N113").Value = -6089

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 35.166668
$ws.Range("I6").Value = 36
$ws.Range("J6").Value = 34.75
$ws.Range("K6").Value = 108
$ws.Range("L6").Value = 104.25
$ws.Range("M6").Value = 5
$ws.Range("N6").Value = -330.25
$ws.Range("H10").Value = 110.5
$ws.Range("I10").Value = 110.5
$ws.Range("K10").Value = 331.5
$ws.Range("M10").Value = -192.5
$ws.Range("H26").Value = 547.5
$ws.Range("I26").Value = 674.375
$ws.Range("J26").Value = 40
$ws.Range("K26").Value = 2023.125
$ws.Range("L26").Value = 120
$ws.Range("M26").Value = -1735.125
$ws.Range("N26").Value = -696

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 703.9524
$ws.Range("I97").Value = 703.9524
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 703.9524
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -207.9524
$ws.Range("N97").ClearContents()
$ws.Range("H132").Value = 142667.38
$ws.Range("I132").Value = 161334.28
$ws.Range("J132").Value = 11999
$ws.Range("K132").Value = 484002.84
$ws.Range("L132").Value = 35997
$ws.Range("M132").Value = -481472.84
$ws.Range("N132").Value = -41057

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H40").Value = 4255.0347
$ws.Range("I40").Value = 2975.087
$ws.Range("J40").Value = 9161.5
$ws.Range("K40").Value = 2975.087
$ws.Range("L40").Value = 9161.5
$ws.Range("M40").Value = -2839.087
$ws.Range("N40").Value = -9433.5
$ws.Range("H68").Value = 9867.166999999999
$ws.Range("J68").Value = 10040.2
$ws.Range("L68").Value = 10040.2
$ws.Range("N68").Value = -11538.2
$ws.Range("H71").Value = 9867.166999999999
$ws.Range("J71").Value = 10040.2
$ws.Range("L71").Value = 50201
$ws.Range("N71").Value = -57689

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 84999.5
$ws.Range("J26").Value = 84999.5
$ws.Range("L26").Value = 84999.5
$ws.Range("N26").Value = -85585.5
$ws.Range("H132").Value = 1703.6744
$ws.Range("I132").Value = 1765
$ws.Range("K132").Value = 5295
$ws.Range("M132").Value = -2765
